# edit.ps1 -- applies the Justification_Report.docx edit described by the
# provided unified diff. The new body content (paragraphs only, matching the
# diff's target state) is embedded below as base64 (UTF-8 XML fragment) to
# sidestep any PowerShell string-quoting/interpolation issues with the
# apostrophes/unicode in the persona text, then grafted onto the document
# in one shot via Range.InsertXML -- the cleanest reliable way to express a
# restructuring this extensive (paragraph merges/splits/insertions) through
# the Word object model.

$d = $word.ActiveDocument

$b64Parts = @(
    "PHc6cD48dzpwUHI+PHc6cFN0eWxlIHc6dmFsPSJIZWFkaW5nMSIvPjwvdzpwUHI+PHc6cj48dzp0PkNTRTJNQUQgUHJvamVjdCBSZXBvcnQgVGVtcGxhdGU8L3c6dD48L3c6cj48dzpyPjx3OmJyLz48dzp0Pkdyb3VwIEluZm9ybWF0aW9uPC93OnQ+PC93OnI+PC93",
    "OnA+PHc6cD48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+R3JvdXAgSWQ6IDwvdzp0PjwvdzpyPjx3OnI+PHc6dD4xNTwvdzp0PjwvdzpyPjwvdzpwPjx3OnA+PHc6cj48dzp0Pk1lbWJlciBOYW1lcyAmYW1wOyBTdHVkZW50IElkczo8L3c6dD48L3c6cj48",
    "L3c6cD48dzpwPjx3OnI+PHc6dD5DYW1lcm9uIFRhbGJvdCAxODUwMTk4NDwvdzp0PjwvdzpyPjwvdzpwPjx3OnA+PC93OnA+PHc6cD48dzpwUHI+PHc6cFN0eWxlIHc6dmFsPSJIZWFkaW5nMSIvPjwvdzpwUHI+PHc6cj48dzp0PkFwcCBvdmVydmlldzwvdzp0Pjwv",
    "dzpyPjwvdzpwPjx3OnA+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPkFwcCB0aGF0IHN0b3JlcyBtdWx0aXBsZSB0eXBlcyBvZiBpbmZvcm1hdGlvbiBhYm91dCBkb3dubG9hZCBzcGVlZHMgd2l0aCBkaWZmZXJlbnQgY29ubmVjdGlvbnMgYW5kIGxvY2F0",
    "aW9ucy4gPC93OnQ+PC93OnI+PC93OnA+PHc6cD48dzpwUHI+PHc6cFN0eWxlIHc6dmFsPSJIZWFkaW5nMSIvPjwvdzpwUHI+PHc6cj48dzp0PlVzZXIgUGVyc29uYS9zPC93OnQ+PC93OnI+PC93OnA+PHc6cD48dzpyPjx3OnQ+RGF2aWQgb3ducyBhIGNhbXAgYW5k",
    "IGNhYmluIGxvZGdlIHdoZXJlIGhlIHJlZ3VsYXJseSBob2xkcyBzcGVjaWFsIGV2ZW50cywgdGVhY2hpbmcgcGVvcGxlIGhvdyB0byBzdXJ2aXZlIGluIG5hdHVyZS4gSGUgaXMgMzQgeWVhcnMgb2xkIGFuZCBsaXZlcyB3aXRoIGhpcyB5b3VuZ2VyIHNpc3RlciBH",
    "d2VuIGluIHRoZSBtYWluIGxvZGdlLjwvdzp0PjwvdzpyPjwvdzpwPjx3OnA+PHc6cj48dzp0PkRhdmlkIGhhcyBiZWVuIHRoaW5raW5nIGFib3V0IGluc3RhbGxpbmcgaW50ZXJuZXQgc29tZXdoZXJlIGluIHRoZSBwYXJrIGFuZCB3b3VsZCBsaWtlIHRvIGtub3cg",
    "d2hhdCBzZXJ2aWNlcyBpbiB0aGUgbG9jYWwgYXJlYSBhcmUgbW9zdCByZWxpYWJsZS4gSGUgaG9wZXMgdGhhdCB0aGUgbmV3IGludGVybmV0IGNvbm5lY3Rpb24gd2lsbCBpbmNyZWFzZSBoaXMgY2xpZW50IGJhc2UgYW5kIG1ha2UgaGlzIHNpc3RlciBoYXBweS48",
    "L3c6dD48L3c6cj48L3c6cD48dzpwPjx3OnBQcj48dzpwQmRyPjx3OmJvdHRvbSB3OnZhbD0ic2luZ2xlIiB3OnN6PSI2IiB3OnNwYWNlPSIxIiB3OmNvbG9yPSJhdXRvIi8+PC93OnBCZHI+PC93OnBQcj48dzpyPjx3OnQ+SGUgd2FudHMgdG8gY29tcGFyZSB0aGUg",
    "Y29uc2lzdGVuY3kgb2Ygc2VydmljZXMgYXQgZGlmZmVyZW50IHRpbWVzIG9mIHRoZSBkYXkgYW5kIHNlZSB3aGljaCBsb2NhdGlvbnMgYXQgdGhlIHBhcmsgaGF2ZSB0aGUgYmVzdCByZWNlcHRpb24uIERhdmlkIGhvd2V2ZXIgaXNu4oCZdCBza2lsbGVkIGluIHVz",
    "aW5nIHRlY2hub2xvZ3kgYW5kIHdvdWxkIGxpa2UgdG8gaGF2ZSBhbiBpbnR1aXRpdmUgYW5kIGVhc3kgdG8gdW5kZXJzdGFuZCBtZXRob2QgdG8gbWFrZSBoaXMgZGVjaXNpb24uPC93OnQ+PC93OnI+PC93OnA+PHc6cD48dzpwUHI+PHc6cEJkcj48dzpib3R0b20g",
    "dzp2YWw9InNpbmdsZSIgdzpzej0iNiIgdzpzcGFjZT0iMSIgdzpjb2xvcj0iYXV0byIvPjwvdzpwQmRyPjwvdzpwUHI+PC93OnA+PHc6cD48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+S2VsbHkgaXMgYSBwcml2YXRlIGludmVzdGlnYXRvci4gSGVyIGpv",
    "YiBpbnZvbHZlcyBjb250YWN0aW5nIGNsaWVudHMgYW5kIGVtYWlsaW5nIHRoZW0gaW5mb3JtYXRpb24gc2hlIGZpbmRzIG9uIGEgcmVndWxhciBiYXNpcyB3aXRoIHVwZGF0ZXMgb24gbmV3IGluZm9ybWF0aW9uIHNoZSBmaW5kcy4gPC93OnQ+PC93OnI+PC93OnA+",
    "PHc6cD48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+U2hlIGlzIDI0IHllYXJzIG9sZCBhbmQgbGl2ZXMgYWxvbmUgaW4gYSBzbWFsbCBhcGFydG1lbnQsIHN0dWR5aW5nIGxhbmd1YWdlIHBhcnQgdGltZS4gU2hlIGhvcGVzIHRvIHRyYXZlbCBhbmQgYWN0",
    "IGFzIGFuIGludGVycHJldGVyIGFuZCB0cmF2ZWwgYWR2aXNvciB0byB0cmF2ZWxsaW5nIGZhbWlsaWVzLiA8L3c6dD48L3c6cj48L3c6cD48dzpwPjx3OnBQcj48dzpwQmRyPjx3OmJvdHRvbSB3OnZhbD0ic2luZ2xlIiB3OnN6PSI2IiB3OnNwYWNlPSIxIiB3OmNv",
    "bG9yPSJhdXRvIi8+PC93OnBCZHI+PC93OnBQcj48dzpyPjx3OnQ+S2VsbHnigJlzIGFwYXJ0bWVudCBpcyBmYWlybHkgc21hbGwgYW5kIHNoZSBsaWtlcyB0byBzcGVuZCBtb3N0IG9mIGhlciB0aW1lIGF0IHRoZSBsb2NhbCBsaWJyYXJ5IGFuZCBjYWbDqXMgd2hp",
    "bGUgc2hlIGRvZXMgaGVyIHJlc2VhcmNoLiBTaGUgd291bGQgbGlrZSB0byBrbm93IHdoZXJlIHRoZSBiZXN0IGludGVybmV0IGlzIHNvIHRoYXQgc2hlIGRvZXNu4oCZdCBoYXZlIHRvIHdhaXQgZm9yIGhlciBzdHVkeSBub3RlcyB0byBsb2FkIGFuZCBzbyBoZXIg",
    "Y2xpZW50cyBjYW4gcmVjZWl2ZSB0aGVpciBlbWFpbHMgYXMgc29vbiBhcyBwb3NzaWJsZS48L3c6dD48L3c6cj48L3c6cD48dzpwPjx3OnBQcj48dzpwQmRyPjx3OmJvdHRvbSB3OnZhbD0ic2luZ2xlIiB3OnN6PSI2IiB3OnNwYWNlPSIxIiB3OmNvbG9yPSJhdXRv",
    "Ii8+PC93OnBCZHI+PC93OnBQcj48L3c6cD48dzpwPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj5MYXJyeSBpcyBhIDwvdzp0PjwvdzpyPjx3OnByb29mRXJyIHc6dHlwZT0iZ3JhbVN0YXJ0Ii8+PHc6cj48dzp0PjE2IHllYXIgb2xkPC93OnQ+PC93OnI+",
    "PHc6cHJvb2ZFcnIgdzp0eXBlPSJncmFtRW5kIi8+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBzdHVkZW50IHdobyBsaXZlcyB3aXRoIGhpcyBtb3RoZXIsIGZhdGhlciBhbmQgNSB5ZWFyIG9sZCBicm90aGVyLiBIZSBpcyBhIHByZXR0eSBhdmVyYWdl",
    "IGtpZCwgZG9lc27igJl0IGdldCBpbnRvIHRyb3VibGUgYXQgc2Nob29sIGFuZCBoYXMgYSBmZXcgZ29vZCBmcmllbmRzIHRoYXQgaGUgc3BlbmRzIHRpbWUgd2l0aCBhZnRlciBzY2hvb2wgYSBmZXcgbmlnaHRzIGEgd2Vlay48L3c6dD48L3c6cj48L3c6cD48dzpw",
    "Pjx3OnI+PHc6dD5XaGVuIExhcnJ5IGlzIGhvbWUgaGUgbGlrZXMgdG8gcGxheSBvbmxpbmUgZ2FtZXMgYW5kIGdldHMgZnJ1c3RyYXRlZCB3aGVuIGhpcyBtdW0gYW5kIGJyb3RoZXIgZ2V0IG9uIHRoZSBpbnRlcm5ldCwgaXQgc2xvd3MgZG93biBoaXMgZ2FtZXMg",
    "YW5kIOKAmG1ha2VzIGhpbSBsb3Nl4oCZLiBMYXJyeSBuZWVkcyB0byBwcm92ZSB0byBoaXMgbXVtIGFuZCBkYWQgdGhhdCB0aGUgaW50ZXJuZXQgc3Vja3Mgc28gdGhhdCB0aGV5IHdpbGwgZ2V0IGEgYmV0dGVyIGludGVybmV0IHBhY2thZ2UuPC93OnQ+PC93OnI+",
    "PC93OnA+PHc6cD48dzpyPjx3Omxhc3RSZW5kZXJlZFBhZ2VCcmVhay8+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj5IaXMgcGFyZW50cyBkb27igJl0IHJlYWxseSB1bmRlcnN0YW5kIHRlY2huaWNhbCBqYXJnb24gc28gaGUgd291bGQgbGlrZSB0byBoYXZlIHJl",
    "c3VsdHMgdGhhdCBhcmUgZWFzeSB0byByZWFkLiBIZSBoYXMgYmVlbiB1c2luZyBhbiBvbmxpbmUgdGVzdGVyIGJ1dCB0aGVyZSBpc27igJl0IGFueXRoaW5nIHRvIGNvbXBhcmUgaGlzIHJlc3VsdHMgdG8uIExhcnJ5IG5lZWRzIHlvdXIgaGVscC4gPC93OnQ+PC93",
    "OnI+PC93OnA+PHc6cD48L3c6cD48dzpwPjx3OnBQcj48dzpwU3R5bGUgdzp2YWw9IkhlYWRpbmcxIi8+PC93OnBQcj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+RXh0ZXJuYWwgbGlicmFyaWVzOiBJZiB5b3UgaGF2ZSB1c2VkIGFueSBleHRlcm5hbCBs",
    "aWJyYXJpZXMgc3VjaCBhcyBzcGVlZC10ZXN0LCA8L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+c3RhdGUgd2hhdCA8L3c6dD48L3c6cj48dzpwcm9vZkVyciB3OnR5cGU9InNwZWxsU3RhcnQiLz48dzpyPjx3OnQ+dGhlPC93OnQ+PC93",
    "OnI+PHc6cHJvb2ZFcnIgdzp0eXBlPSJzcGVsbEVuZCIvPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gPC93OnQ+PC93OnI+PHc6cHJvb2ZFcnIgdzp0eXBlPSJzcGVsbFN0YXJ0Ii8+PHc6cj48dzp0PmFyZTwvdzp0PjwvdzpyPjx3OnByb29mRXJyIHc6",
    "dHlwZT0ic3BlbGxFbmQiLz48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IGFuZCA8L3c6dD48L3c6cj48dzpyPjx3OnQ+anVzdGlmeSByZWFzb25zIGZvciBkb2luZyBzbzwvdzp0PjwvdzpyPjwvdzpwPjx3OnA+PHc6cFByPjx3OnBTdHlsZSB3OnZhbD0i",
    "SGVhZGluZzEiLz48L3c6cFByPjwvdzpwPjx3OnA+PHc6cFByPjx3OnBTdHlsZSB3OnZhbD0iSGVhZGluZzEiLz48L3c6cFByPjx3OnI+PHc6dD5UZWFtIE1hbmFnZW1lbnQ8L3c6dD48L3c6cj48L3c6cD48dzpwPjx3OnBQcj48dzpwU3R5bGUgdzp2YWw9IkhlYWRp",
    "bmcyIi8+PHc6bnVtUHI+PHc6aWx2bCB3OnZhbD0iMCIvPjx3Om51bUlkIHc6dmFsPSIxIi8+PC93Om51bVByPjx3OnJQcj48dzpyRm9udHMgdzphc2NpaT0iVGltZXMiIHc6ZWFzdEFzaWE9IlRpbWVzIE5ldyBSb21hbiIgdzpoQW5zaT0iVGltZXMiIHc6Y3M9IlRp",
    "bWVzIE5ldyBSb21hbiIvPjx3OnN6IHc6dmFsPSIyMCIvPjx3OnN6Q3Mgdzp2YWw9IjIwIi8+PC93OnJQcj48L3c6cFByPjx3OmJvb2ttYXJrU3RhcnQgdzppZD0iMCIgdzpuYW1lPSJfR29CYWNrIi8+PHc6Ym9va21hcmtFbmQgdzppZD0iMCIvPjx3OnI+PHc6clBy",
    "Pjx3OnJGb250cyB3OmFzY2lpPSJUaW1lcyIgdzplYXN0QXNpYT0iVGltZXMgTmV3IFJvbWFuIiB3OmhBbnNpPSJUaW1lcyIgdzpjcz0iVGltZXMgTmV3IFJvbWFuIi8+PHc6c3ogdzp2YWw9IjIwIi8+PHc6c3pDcyB3OnZhbD0iMjAiLz48L3c6clByPjx3OnQ+QnVy",
    "bmRvd24gY2hhcnQgZm9yIFNwcmludCAyPC93OnQ+PC93OnI+PC93OnA+PHc6cD48dzpwUHI+PHc6cFN0eWxlIHc6dmFsPSJIZWFkaW5nMiIvPjx3Om51bVByPjx3Omlsdmwgdzp2YWw9IjAiLz48dzpudW1JZCB3OnZhbD0iMSIvPjwvdzpudW1Qcj48dzpyUHI+PHc6",
    "ckZvbnRzIHc6YXNjaWk9IlRpbWVzIiB3OmVhc3RBc2lhPSJUaW1lcyBOZXcgUm9tYW4iIHc6aEFuc2k9IlRpbWVzIiB3OmNzPSJUaW1lcyBOZXcgUm9tYW4iLz48dzpzeiB3OnZhbD0iMjAiLz48dzpzekNzIHc6dmFsPSIyMCIvPjwvdzpyUHI+PC93OnBQcj48dzpy",
    "Pjx3OnJQcj48dzpyRm9udHMgdzphc2NpaT0iVGltZXMiIHc6ZWFzdEFzaWE9IlRpbWVzIE5ldyBSb21hbiIgdzpoQW5zaT0iVGltZXMiIHc6Y3M9IlRpbWVzIE5ldyBSb21hbiIvPjx3OnN6IHc6dmFsPSIyMCIvPjx3OnN6Q3Mgdzp2YWw9IjIwIi8+PC93OnJQcj48",
    "dzp0PlNwcmludCBiYWNrbG9nIGZvciBuZXh0IHNwcmludCAoU3ByaW50IDMpPC93OnQ+PC93OnI+PC93OnA+PHc6cD48dzpwUHI+PHc6cFN0eWxlIHc6dmFsPSJIZWFkaW5nMiIvPjwvdzpwUHI+PC93OnA+"
)
$b64 = [string]::Join("", $b64Parts)
$bytes = [System.Convert]::FromBase64String($b64)
$bodyXml = [System.Text.Encoding]::UTF8.GetString($bytes)

$pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  $bodyXml +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Graft the new paragraph sequence over the whole body (sectPr at the end is
# untouched since $d.Content stops just before it).
$d.Content.InsertXML($pkgXml)

# The cached PAGE field result in the default footer bumps from 1 to 2 now
# that the report spans an extra page.
$sec = $d.Sections(1)
$footer = $sec.Footers(1)
$null = $footer.Range.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "2", 2)

Write-Host "Paragraphs:" $d.Paragraphs.Count
